# "Add files via upload" - fill in the script-name column (A) for the
# newly-added rows in the transition-exit table, and move the selection /
# scrolled viewport down to the rows that were just filled in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 reuses the script name already used two rows above (row 19 / D73P29A)
$ws.Range("A20").Value = "SCRIPT/D73P29A/enter06.ssb"

# Rows 21, 22 and 48 introduce brand-new script names (added to the shared
# string table in this order), row 23 reuses the last of the four new ones.
$ws.Range("A21").Value = "SCRIPT/D79P21A/enter12.ssb"
$ws.Range("A22").Value = "SCRIPT/D35P21A/enter12.ssb"
$ws.Range("A48").Value = "SCRIPT/D39P21A/enter12.ssb"
$ws.Range("A23").Value = "SCRIPT/D41P21A/enter12.ssb"

# Scroll the window down a bit and move the active selection from E48 to
# D48, matching where the author ended up after entering the new data.
$win = $excel.ActiveWindow
$win.ScrollRow = 46
$win.ScrollColumn = 1
$ws.Range("D48").Select()
